# Generate Report for Handoff
# Adds two new localization entries
#   e842a4a6-e166-497d-acae-045c424a470f.md
#   e88d6b61-314d-4175-9ac7-1d8e633c9cad.md
# as new rows to the Overview / zh-cn / de-de tables.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"
$hlColor = 15570276   # RGB(100,149,237) == FF6495ED, matches the workbook's HyperLink style
$hlUnderline = 2      # xlUnderlineStyleSingle

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = $hlUnderline
    $rng.Font.Color = $hlColor
}

# ----------------------------------------------------------------------------
# Sheet "Overview"
# ----------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

# Row 4
$wsOverview.Range("A4").Value2 = "e842a4a6-e166-497d-acae-045c424a470f.md"
$wsOverview.Range("B4").Value2 = "e2e\e842a4a6-e166-497d-acae-045c424a470f.md"
Style-AsHyperlink($wsOverview.Range("B4"))
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4706de31538cd1d3573b3b0fa014e3bf32f343a9/e2e/e842a4a6-e166-497d-acae-045c424a470f.md", "", "", "e2e\e842a4a6-e166-497d-acae-045c424a470f.md") | Out-Null
$wsOverview.Range("C4").Value2 = ".md"
$wsOverview.Range("D4").Value2 = ""
$wsOverview.Range("E4").Value2 = "Ready for handoff"
$wsOverview.Range("F4").Value2 = "Ready for handoff"
$wsOverview.Range("G4").Value2 = "2016-08-25 02:40:26"
$wsOverview.Range("G4").NumberFormat = $dateFmt

# Row 5
$wsOverview.Range("A5").Value2 = "e88d6b61-314d-4175-9ac7-1d8e633c9cad.md"
$wsOverview.Range("B5").Value2 = "e2e\e88d6b61-314d-4175-9ac7-1d8e633c9cad.md"
Style-AsHyperlink($wsOverview.Range("B5"))
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40b016d8fd6344d61935726d5ee2194d0166d365/e2e/e88d6b61-314d-4175-9ac7-1d8e633c9cad.md", "", "", "e2e\e88d6b61-314d-4175-9ac7-1d8e633c9cad.md") | Out-Null
$wsOverview.Range("C5").Value2 = ".md"
$wsOverview.Range("D5").Value2 = ""
$wsOverview.Range("E5").Value2 = "Ready for handoff"
$wsOverview.Range("F5").Value2 = "Ready for handoff"
$wsOverview.Range("G5").Value2 = "2016-08-25 02:40:26"
$wsOverview.Range("G5").NumberFormat = $dateFmt

# ----------------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)

$loZhCn.ListRows.Add() | Out-Null
$loZhCn.ListRows.Add() | Out-Null

# Row 4
$wsZhCn.Range("A4").Value2 = "e842a4a6-e166-497d-acae-045c424a470f.md"
Style-AsHyperlink($wsZhCn.Range("A4"))
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4706de31538cd1d3573b3b0fa014e3bf32f343a9/e2e/e842a4a6-e166-497d-acae-045c424a470f.md", "", "", "e842a4a6-e166-497d-acae-045c424a470f.md") | Out-Null
$wsZhCn.Range("B4").Value2 = ".md"
$wsZhCn.Range("C4").Value2 = "Ready for handoff"
$wsZhCn.Range("D4").Value2 = "e2e"
$wsZhCn.Range("E4").Value2 = "ht"
$wsZhCn.Range("F4").Value2 = "'False"
$wsZhCn.Range("G4").Value2 = "e842a4a6-e166-497d-acae-045c424a470f.4706de31538cd1d3573b3b0fa014e3bf32f343a9.zh-cn.xlf"
$wsZhCn.Range("H4").Value2 = "2016-08-25 02:40:21"
$wsZhCn.Range("H4").NumberFormat = $dateFmt
$wsZhCn.Range("I4").Value2 = ""
$wsZhCn.Range("J4").Value2 = ""
$wsZhCn.Range("K4").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("K4").NumberFormat = $dateFmt
$wsZhCn.Range("L4").Value2 = ""
$wsZhCn.Range("M4").Value2 = "'True"
$wsZhCn.Range("N4").Value2 = ""
$wsZhCn.Range("O4").Value2 = "'False"
$wsZhCn.Range("P4").Value2 = ""

# Row 5
$wsZhCn.Range("A5").Value2 = "e88d6b61-314d-4175-9ac7-1d8e633c9cad.md"
Style-AsHyperlink($wsZhCn.Range("A5"))
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40b016d8fd6344d61935726d5ee2194d0166d365/e2e/e88d6b61-314d-4175-9ac7-1d8e633c9cad.md", "", "", "e88d6b61-314d-4175-9ac7-1d8e633c9cad.md") | Out-Null
$wsZhCn.Range("B5").Value2 = ".md"
$wsZhCn.Range("C5").Value2 = "Ready for handoff"
$wsZhCn.Range("D5").Value2 = "e2e"
$wsZhCn.Range("E5").Value2 = "ht"
$wsZhCn.Range("F5").Value2 = "'False"
$wsZhCn.Range("G5").Value2 = "e88d6b61-314d-4175-9ac7-1d8e633c9cad.40b016d8fd6344d61935726d5ee2194d0166d365.zh-cn.xlf"
$wsZhCn.Range("H5").Value2 = "2016-08-25 02:40:21"
$wsZhCn.Range("H5").NumberFormat = $dateFmt
$wsZhCn.Range("I5").Value2 = ""
$wsZhCn.Range("J5").Value2 = ""
$wsZhCn.Range("K5").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("K5").NumberFormat = $dateFmt
$wsZhCn.Range("L5").Value2 = ""
$wsZhCn.Range("M5").Value2 = "'True"
$wsZhCn.Range("N5").Value2 = ""
$wsZhCn.Range("O5").Value2 = "'False"
$wsZhCn.Range("P5").Value2 = ""

# ----------------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)

$loDeDe.ListRows.Add() | Out-Null
$loDeDe.ListRows.Add() | Out-Null

# Row 4
$wsDeDe.Range("A4").Value2 = "e842a4a6-e166-497d-acae-045c424a470f.md"
Style-AsHyperlink($wsDeDe.Range("A4"))
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4706de31538cd1d3573b3b0fa014e3bf32f343a9/e2e/e842a4a6-e166-497d-acae-045c424a470f.md", "", "", "e842a4a6-e166-497d-acae-045c424a470f.md") | Out-Null
$wsDeDe.Range("B4").Value2 = ".md"
$wsDeDe.Range("C4").Value2 = "Ready for handoff"
$wsDeDe.Range("D4").Value2 = "e2e"
$wsDeDe.Range("E4").Value2 = "ht"
$wsDeDe.Range("F4").Value2 = "'False"
$wsDeDe.Range("G4").Value2 = "e842a4a6-e166-497d-acae-045c424a470f.4706de31538cd1d3573b3b0fa014e3bf32f343a9.de-de.xlf"
$wsDeDe.Range("H4").Value2 = "2016-08-25 02:40:26"
$wsDeDe.Range("H4").NumberFormat = $dateFmt
$wsDeDe.Range("I4").Value2 = ""
$wsDeDe.Range("J4").Value2 = ""
$wsDeDe.Range("K4").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("K4").NumberFormat = $dateFmt
$wsDeDe.Range("L4").Value2 = ""
$wsDeDe.Range("M4").Value2 = "'True"
$wsDeDe.Range("N4").Value2 = ""
$wsDeDe.Range("O4").Value2 = "'False"
$wsDeDe.Range("P4").Value2 = ""

# Row 5
$wsDeDe.Range("A5").Value2 = "e88d6b61-314d-4175-9ac7-1d8e633c9cad.md"
Style-AsHyperlink($wsDeDe.Range("A5"))
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40b016d8fd6344d61935726d5ee2194d0166d365/e2e/e88d6b61-314d-4175-9ac7-1d8e633c9cad.md", "", "", "e88d6b61-314d-4175-9ac7-1d8e633c9cad.md") | Out-Null
$wsDeDe.Range("B5").Value2 = ".md"
$wsDeDe.Range("C5").Value2 = "Ready for handoff"
$wsDeDe.Range("D5").Value2 = "e2e"
$wsDeDe.Range("E5").Value2 = "ht"
$wsDeDe.Range("F5").Value2 = "'False"
$wsDeDe.Range("G5").Value2 = "e88d6b61-314d-4175-9ac7-1d8e633c9cad.40b016d8fd6344d61935726d5ee2194d0166d365.de-de.xlf"
$wsDeDe.Range("H5").Value2 = "2016-08-25 02:40:26"
$wsDeDe.Range("H5").NumberFormat = $dateFmt
$wsDeDe.Range("I5").Value2 = ""
$wsDeDe.Range("J5").Value2 = ""
$wsDeDe.Range("K5").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("K5").NumberFormat = $dateFmt
$wsDeDe.Range("L5").Value2 = ""
$wsDeDe.Range("M5").Value2 = "'True"
$wsDeDe.Range("N5").Value2 = ""
$wsDeDe.Range("O5").Value2 = "'False"
$wsDeDe.Range("P5").Value2 = ""

Write-Host "Report generated for handoff."
